$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, copying the formatting from the existing
# header cell H1 (bold, centered, bordered) so the new headers match style.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the I and J columns with their numeric values for rows 2-31.
$iValues = @(7,8,8,8,5,6,6,8,8,7,11,6,7,9,8,7,7,8,6,8,7,6,6,6,7,3,7,9,8,2)
$jValues = @(9,8,8,8,6,6,6,8,8,7,11,6,8,9,8,7,7,8,6,8,7,6,6,6,7,4,8,9,8,2)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
